# UC003 - Cancelar Solicitação de Diária: version bump + text fixes + TC3/TC4 step restructuring
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -------------------------------------------------------------------
# 1. Version: 1.0 -> 1.2.5
# -------------------------------------------------------------------
$ws.Range("D2").Value = "1.2.5"

# -------------------------------------------------------------------
# 2. Precondition text fix (shared across all 5 test-case blocks):
#    "O usuario ..." -> "O usuário ... diárias." (accent + trailing period)
# -------------------------------------------------------------------
$precondition = "O usuário devidamente autenticado e na tela inicial de cancelar diárias."
$ws.Range("B8").Value = $precondition
$ws.Range("B17").Value = $precondition
$ws.Range("B25").Value = $precondition
$ws.Range("B32").Value = $precondition
$ws.Range("B40").Value = $precondition

# -------------------------------------------------------------------
# 3. MSG102 text fix (trailing period added), used by TC1/TC2/TC5 steps
# -------------------------------------------------------------------
$msg102 = "SYSTEM Exibe a mensagem (MSG102 - Confirmar cancelamento)."
$ws.Range("D10").Value = $msg102
$ws.Range("D19").Value = $msg102
$ws.Range("D42").Value = $msg102

# -------------------------------------------------------------------
# 4. MSG217 text fix (stray tab character removed), used by TC2 step 2
# -------------------------------------------------------------------
$msg217 = "SYSTEM Identifica que o usuário não informou uma justificativa para o cancelamento. Não efetiva o cancelamento e exibe mensagem de erro (MSG217 - Necessário informar uma justificativa para o cancelamento de solicitações) para o usuário."
$ws.Range("D20").Value = $msg217

# -------------------------------------------------------------------
# 5. Restructure TC3 (now 2 steps, same as TC2 pattern) and TC4
#    (now 1 step, with MSG205 typo fixed). TC3 currently has 1 step
#    row (27) followed by two blank rows (28,29) before TC4's header
#    (30). TC4 currently has 2 step rows (34,35).
#
#    Net effect needed: insert a row before the current blank row 28
#    (giving TC3 a second step row) and delete the row that held
#    TC4's second step (so TC5, starting at row 38, does not move).
# -------------------------------------------------------------------

# Insert a new row at 28 - this shifts old rows 28-44 down to 29-45.
$ws.Rows("28:28").Insert()

# Copy the formatting (borders/fonts/number format) of the TC3 step-1
# row (27) onto the freshly inserted row 28, so it matches the visual
# style of a normal step row.
$ws.Range("A27:F27").Copy()
$ws.Range("A28:F28").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# TC3 step 1 (row 27): "Chefe Não informa..." / "SYSTEM Exibe MSG102..."
$ws.Range("A27").Value = 1.0
$ws.Range("B27").Value = "Chefe Não informa o motivo do cancelamento."
$ws.Range("D27").Value = $msg102

# TC3 step 2 (row 28, new): "Chefe Clica em confirmar." / MSG217 text
$ws.Range("A28").Value = 2.0
$ws.Range("B28").Value = "Chefe Clica em confirmar."
$ws.Range("D28").Value = $msg217

# After the insert, TC4's block (old rows 30-35) now lives at rows
# 31-36: header(31) desc(32) precondition(33) steps-header(34)
# step1(35, old content) step2(36, old content - to be removed).

# TC4 step 1 (row 35, after shift): "Chefe Informa..." / MSG205 (typo fixed)
$msg205 = "SYSTEM Identifica que a solicitação de diária está em situação diferente de 'SOLICITADA PARA EMPENHO' ou 'SOLICITADA PARA PRESTAÇÃO DE CONTAS'.  Impede o cancelamento e exibe mensagem de erro (MSG205 - Solicitação de diária não pode ser cancelada) para o usuário."
$ws.Range("A35").Value = 1.0
$ws.Range("B35").Value = "Chefe Informa o motivo do cancelamento."
$ws.Range("D35").Value = $msg205

# Remove the now-obsolete second TC4 step row (old step 2, currently
# at row 36) so TC5's header stays anchored at row 38.
$ws.Rows("36:36").Delete()
